# Swap the weekly price-report values between row pairs (2,6), (3,7), (4,8), (5,9).
# Only columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) are affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "M", "N", "O", "P", "S")
$pairs = @(
    @(2, 6),
    @(3, 7),
    @(4, 8),
    @(5, 9)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
